# Weekly update: insert a new price record (row 18) for "Feria Lagunitas de
# Puerto Montt" - Alcachofa, pushing the existing historical rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18; this shifts rows 18..44 down to 19..45
# and keeps the "Fecha" column's date number format (style index 2) intact.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(18, 1).Value  = 4
$ws.Cells.Item(18, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(18, 3).Value  = "Los Lagos"
$ws.Cells.Item(18, 4).Value  = 45195
$ws.Cells.Item(18, 5).Value  = 10
$ws.Cells.Item(18, 6).Value  = 100112013
$ws.Cells.Item(18, 7).Value  = "Alcachofa"
$ws.Cells.Item(18, 8).Value  = "Argentina(o)"
$ws.Cells.Item(18, 9).Value  = "Primera"
$ws.Cells.Item(18, 10).Value = 150
$ws.Cells.Item(18, 11).Value = 14000
$ws.Cells.Item(18, 12).Value = 14000
$ws.Cells.Item(18, 13).Value = 14000
$ws.Cells.Item(18, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(18, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 16).Value = 280
$ws.Cells.Item(18, 17).Value = 50
$ws.Cells.Item(18, 18).Value = "Hortaliza"
